# Update "想去人数" (interested-count) figures in the "展览" (sheet 1)
# and "全部类型" (sheet 4) worksheets, matching the freshly generated
# gh-pages data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item(1)   # 展览
$sheetAllTypes    = $wb.Worksheets.Item(4)  # 全部类型

# Row -> new F-column value, for the "展览" sheet (rows 2-24)
$exhibitionUpdates = @{
    2  = 1577
    3  = 8929
    5  = 501
    6  = 677
    7  = 335
    9  = 41
    10 = 59
    11 = 3789
    13 = 376
    15 = 4144
    18 = 1137
    19 = 3
    23 = 2591
    24 = 100
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F-column value, for the "全部类型" sheet (rows 2-25)
$allTypesUpdates = @{
    2  = 1577
    3  = 8929
    5  = 501
    6  = 677
    7  = 335
    9  = 41
    10 = 59
    11 = 3789
    13 = 376
    15 = 4144
    18 = 1137
    19 = 3
    23 = 2591
    25 = 100
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
